$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# EPBDS-14389 Implement tryJSON() method
#
# Appends a new "SimpleRules String validate(State st)" table (rows 41-48,
# columns A-E) below the existing "Datatype Complex" table (rows 36-40),
# mirroring the existing sheet's formatting conventions.
# ---------------------------------------------------------------------------

# Stamp the default (gray, bordered, blank) style across the whole new block
# first, by copying the format of an existing blank separator row (row 36).
# This reuses the workbook's existing style index instead of synthesizing a
# brand-new one.
$ws.Range("A36:E36").Copy()
$ws.Range("A41:E48").PasteSpecial(-4122)

# Restore the explicit custom row height used throughout the sheet for every
# newly added row.
$ws.Rows.Item(41).RowHeight() = 13.55
$ws.Rows.Item(42).RowHeight() = 13.55
$ws.Rows.Item(43).RowHeight() = 13.55
$ws.Rows.Item(44).RowHeight() = 13.55
$ws.Rows.Item(45).RowHeight() = 13.55
$ws.Rows.Item(46).RowHeight() = 13.55
$ws.Rows.Item(47).RowHeight() = 13.55
$ws.Rows.Item(48).RowHeight() = 13.55

# Table header.
$ws.Range("B42").Value() = "SimpleRules String validate(State st)"

# Column headers.
$ws.Range("B43").Value() = "State"
$ws.Range("C43").Value() = "Result"

# Rule rows - plain text values. Set left-to-right, top-to-bottom so new
# shared-string indices are allocated in the same order as the source table.
$ws.Range("B44").Value() = "NY"
$ws.Range("C44").Value() = "OK"

$ws.Range("B45").Value() = "CA"
# This result cell is literal text starting with "=" (OpenL rule syntax),
# not a real formula - a leading apostrophe forces Excel to store it as
# plain text instead of parsing/evaluating it.
$ws.Range("C45").Value() = "'= error(""CA is not allowed"")"

$ws.Range("B46").Value() = "MI"
$ws.Range("C46").Value() = "'= error(""CD1"", ""Failure"")"

$ws.Range("C47").Value() = "'= error(new Complex(""Yura"", 1000))"

# Re-apply the shaded/text style (index 13, same as the existing
# "Datatype Complex" table body) to every cell that now holds a shared
# string - setting a quote-prefixed value above mutates the cell's style
# (adds a transient quotePrefix flag), so this pass normalizes them all
# back to the canonical style used elsewhere on the sheet.
$ws.Range("B37:B39").Copy()
$ws.Range("B42:B46").PasteSpecial(-4122)

$ws.Range("C38:C39").Copy()
$ws.Range("C43:C47").PasteSpecial(-4122)

Write-Output "done"
